# Add a new "staff" column (K) to the sample data sheet, with a value
# for each of the three hospital rows, matching the header/number
# formatting already used for column J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell — copy J1's format (style) onto K1, then set its text.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "staff"

# Data cells — copy each row's J-column number format onto K, then set
# the new staff-count values.
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("K2").Value = 275

$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = 425

$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 394
